# excel data driven integration
# Update the sample email addresses used in the "registerUsr" test data sheet
# and leave the active selection on the Email column for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registerUsr")

$ws.Range("G2").Value = "vijeysssdaayssww@gmail.com"
$ws.Range("G3").Value = "divyammmqqyssww@gmail.com"

$ws.Activate()
$ws.Range("G8").Select()
